$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($row, $col, $value)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextCell 2 4 "29.960.43"
Set-TextCell 2 5 "  +0.54%  "
Set-TextCell 3 4 "1.892.64"
Set-TextCell 3 5 "  +0.06%  "
Set-TextCell 4 5 "  -0.06%  "
Set-TextCell 5 4 "0.7741"
Set-TextCell 5 5 "  -0.72%  "
Set-TextCell 6 4 "243.82"
Set-TextCell 6 5 "  +0.03%  "
Set-TextCell 7 5 "  -0.05%  "
Set-TextCell 8 4 "0.3130"
Set-TextCell 8 5 "  +0.17%  "
Set-TextCell 9 4 "25.79"
Set-TextCell 9 5 "  +2.31%  "
Set-TextCell 10 4 "0.07261"
Set-TextCell 11 4 "0.08702"
Set-TextCell 11 5 "  +7.93%  "
Set-TextCell 12 4 "1.992.04"
Set-TextCell 12 5 "  +5.71%  "
Set-TextCell 13 4 "0.7731"
Set-TextCell 13 5 "  +1.59%  "
Set-TextCell 14 5 "  -0.41%  "
Set-TextCell 15 4 "94.55"
Set-TextCell 15 5 "  +2.70%  "
Set-TextCell 16 4 "6.210"
Set-TextCell 16 5 "  +1.10%  "
Set-TextCell 17 4 "30.076.84"
Set-TextCell 17 5 "  +0.88%  "
Set-TextCell 18 5 "  +0.27%  "
Set-TextCell 19 4 "245.64"
Set-TextCell 19 5 "  +1.05%  "
Set-TextCell 20 4 "2.294.87"
Set-TextCell 20 5 "  +6.65%  "
Set-TextCell 21 5 "  +1.70%  "
Set-TextCell 22 4 "8.191"
Set-TextCell 22 5 "  +1.31%  "
Set-TextCell 23 4 "1.001"
Set-TextCell 23 5 "  -0.02%  "
Set-TextCell 24 4 "1.001"
Set-TextCell 24 5 "  -0.06%  "
Set-TextCell 25 4 "0.1601"
Set-TextCell 25 5 "  -1.47%  "
Set-TextCell 26 4 "9.540"
Set-TextCell 26 5 "  +1.78%  "
Set-TextCell 27 4 "162.96"
Set-TextCell 27 5 "  +0.09%  "
Set-TextCell 28 4 "18.85"
Set-TextCell 28 5 "  +0.95%  "
Set-TextCell 29 4 "2.048"
Set-TextCell 29 5 "  +0.30%  "
Set-TextCell 30 4 "1.431"
Set-TextCell 30 5 "  +1.70%  "
Set-TextCell 31 4 "1.546"
Set-TextCell 31 5 "  +0.06%  "
Set-TextCell 32 4 "4.537"
Set-TextCell 32 5 "  +1.61%  "
Set-TextCell 33 4 "4.131"
Set-TextCell 33 5 "  +1.03%  "
Set-TextCell 34 4 "0.05447"
Set-TextCell 34 5 "  -1.14%  "
Set-TextCell 35 4 "1.250"
Set-TextCell 35 5 "  -1.03%  "
Set-TextCell 36 4 "0.7547"
Set-TextCell 36 5 "  +1.75%  "
Set-TextCell 37 4 "1.000"
Set-TextCell 37 5 "  +0.29%  "
Set-TextCell 38 4 "2.688"
Set-TextCell 38 5 "  +2.45%  "
Set-TextCell 39 4 "0.01973"
Set-TextCell 39 5 "  +3.13%  "
Set-TextCell 40 4 "2.784"
Set-TextCell 40 5 "  +0.35%  "
Set-TextCell 41 4 "0.4525"
Set-TextCell 41 5 "  +2.71%  "
Set-TextCell 42 4 "73.84"
Set-TextCell 42 5 "  +0.58%  "
Set-TextCell 43 2 "Maker"
Set-TextCell 43 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell 43 4 "1.094.87"
Set-TextCell 43 5 "  -3.74%  "
Set-TextCell 44 2 "FraxShare"
Set-TextCell 44 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell 44 4 "6.064"
Set-TextCell 44 5 "  +3.95%  "
Set-TextCell 45 5 "  +0.04%  "
Set-TextCell 46 2 "PaxDollar"
Set-TextCell 46 3 "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell 46 4 "1.000"
Set-TextCell 46 5 "  -0.07%  "
Set-TextCell 47 2 "RocketPoolETH"
Set-TextCell 47 3 "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell 47 4 "2.191.31"
Set-TextCell 47 5 "  +7.22%  "
Set-TextCell 48 4 "103.34"
Set-TextCell 48 5 "  -0.19%  "
Set-TextCell 49 5 "  +1.15%  "
Set-TextCell 50 4 "7.629"
Set-TextCell 50 5 "  +2.82%  "
Set-TextCell 51 4 "9.859"
Set-TextCell 51 5 "  -0.54%  "
